$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the long multi-step procedure texts in column D with short
# "<Action> Setup Subsektor" labels.
$ws.Range("D2").Value = "Tambah Setup Subsektor"
$ws.Range("D3").Value = "View Setup Subsektor"
$ws.Range("D4").Value = "Ubah Setup Subsektor"
$ws.Range("D5").Value = "Hapus Setup Subsektor"

# Row heights shrink now that the cell text is much shorter (wrap text
# still applies via the existing cell style).
$ws.Rows(2).RowHeight = 30
$ws.Rows(3).RowHeight = 30
$ws.Rows(4).AutoFit()
$ws.Rows(5).RowHeight = 30

# Move the sheet's active selection to D5.
$ws.Range("D5").Select()
